$d = $word.ActiveDocument

$bullet = [char]0x2022
$pm = [char]0xB1
$metricColor = 5258796   # RGB(0x2C,0x3E,0x50) packed as a BGR long for Word's Font.Color

function Get-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

function Set-MetricBoldColor($para, $searchText) {
    $r = $para.Range
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "MISSING MATCH: $searchText"
        return
    }
    $r.Font.Bold = 1
    $r.Font.Color = $metricColor
}

# 1) Partner - Siege Analytics achievement bullet: 23% / 64%
$text1 = "$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%"
$p1 = Get-ParagraphByText $d $text1
if ($p1 -ne $null) {
    Set-MetricBoldColor $p1 "23%"
    Set-MetricBoldColor $p1 "64%"
} else {
    Write-Output "Paragraph 1 not found"
}

# 2) Partner - Siege Analytics achievement bullet: 87% / 71% / ±4.2% / ±2.1%
$text2 = "$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${pm}4.2% to ${pm}2.1%"
$p2 = Get-ParagraphByText $d $text2
if ($p2 -ne $null) {
    Set-MetricBoldColor $p2 "87%"
    Set-MetricBoldColor $p2 "71%"
    $margin1 = "${pm}4.2%"
    $margin2 = "${pm}2.1%"
    Set-MetricBoldColor $p2 $margin1
    Set-MetricBoldColor $p2 $margin2
} else {
    Write-Output "Paragraph 2 not found"
}

# 3) Senior Analyst - Myers Research bullet: 1,200
$text3 = "$bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
$p3 = Get-ParagraphByText $d $text3
if ($p3 -ne $null) {
    Set-MetricBoldColor $p3 "1,200"
} else {
    Write-Output "Paragraph 3 not found"
}

# 4) Lake Research Partners bullet: $400M / $1B
$dollar = [char]0x24
$text4 = "$bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the ${dollar}400M Polling Consortium Database at The Analyst Institute, now valued at ${dollar}1B+"
$p4 = Get-ParagraphByText $d $text4
if ($p4 -ne $null) {
    $m400 = "${dollar}400M"
    $m1b = "${dollar}1B"
    Set-MetricBoldColor $p4 $m400
    Set-MetricBoldColor $p4 $m1b
} else {
    Write-Output "Paragraph 4 not found"
}

# 5) Key Achievements bullet: 73.5% / $4.7M
$text5 = "$bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations ${dollar}4.7M"
$p5 = Get-ParagraphByText $d $text5
if ($p5 -ne $null) {
    $m47 = "${dollar}4.7M"
    Set-MetricBoldColor $p5 "73.5%"
    Set-MetricBoldColor $p5 $m47
} else {
    Write-Output "Paragraph 5 not found"
}

# 6) Key Achievements bullet: 87% / 71% (shorter variant)
$text6 = "$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
$p6 = Get-ParagraphByText $d $text6
if ($p6 -ne $null) {
    Set-MetricBoldColor $p6 "87%"
    Set-MetricBoldColor $p6 "71%"
} else {
    Write-Output "Paragraph 6 not found"
}
